$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1791714497743671
$ws.Range("C2").Value = 0.5385450383678798
$ws.Range("D2").Value = 0.6945651025515058
$ws.Range("E2").Value = 0.8334057250532335
$ws.Range("F2").Value = 0.8446427427717865

$ws.Range("B3").Value = 0.301941271613995
$ws.Range("C3").Value = 0.4442637132509937
$ws.Range("D3").Value = 0.4245434332406184
$ws.Range("E3").Value = 0.6515699757053101
$ws.Range("F3").Value = 0.600962680661734

$ws.Range("B4").Value = 0.3243108425389511
$ws.Range("C4").Value = 0.4322851197510529
$ws.Range("D4").Value = 0.3306362357032265
$ws.Range("E4").Value = 0.5750097700937146
$ws.Range("F4").Value = 0.495938463684469

$ws.Range("B5").Value = 0.3104853196909608
$ws.Range("C5").Value = 0.3658015102623823
$ws.Range("D5").Value = 0.2953736427042636
$ws.Range("E5").Value = 0.5434828817030611
$ws.Range("F5").Value = 0.4678351845006229

$ws.Range("B6").Value = 0.268364776153357
$ws.Range("C6").Value = 0.3511917432508606
$ws.Range("D6").Value = 0.2725367241778094
$ws.Range("E6").Value = 0.5220504996432906
$ws.Range("F6").Value = 0.4720135015700387

$ws.Range("B7").Value = 0.262150489134635
$ws.Range("C7").Value = 0.3806290492689116
$ws.Range("D7").Value = 0.3627195589636479
$ws.Range("E7").Value = 0.6022620351339173
$ws.Range("F7").Value = 0.5751054381688496

$ws.Range("B8").Value = 0.4100380287874318
$ws.Range("C8").Value = 0.5602907672914271
$ws.Range("D8").Value = 0.4907109460630159
$ws.Range("E8").Value = 0.7005076345501282
$ws.Range("F8").Value = 0.6221701641941373

$ws.Range("B9").Value = 0.7531841389519395
$ws.Range("C9").Value = 0.7531841389519395
$ws.Range("D9").Value = 0.7596698399655647
$ws.Range("E9").Value = 0.8715904083717103
$ws.Range("F9").Value = 0.5371919947236602

$ws.Range("B10").Value = -0.07569287514151168
$ws.Range("C10").Value = 0.07569287514151168
$ws.Range("D10").Value = 0.005729411347188478
$ws.Range("E10").Value = 0.07569287514151168
